$d = $word.ActiveDocument

# The document body ends with two empty paragraphs followed by the
# section properties. We need to insert a brand-new paragraph containing
# the "Opmerking:" note between those two empty paragraphs, leaving both
# empty paragraphs untouched.

$paraCount = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($paraCount)

# Insert a fresh empty paragraph right before the final (second) empty
# paragraph - this lands it directly after the first empty paragraph.
$lastPara.Range.InsertParagraphBefore()

# The paragraph we just created is now the second-to-last paragraph.
$newParaCount = $d.Paragraphs.Count
$newPara = $d.Paragraphs.Item($newParaCount - 1)
$newRange = $newPara.Range

# Build the new paragraph's content as two runs (matching how Word splits
# runs at a manual line break): one run with "Opmerking:" and a second
# run starting with the line break followed by the note text.
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Opmerking:</w:t></w:r><w:r><w:br/><w:t>De timestamp geeft de verkeerde datum terug, geeft een tijd terug op 1 januari 2001</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$newRange.InsertXML($xml)
